$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 364.45456
$ws.Range("I33").Value = 360.9
$ws.Range("K33").Value = 360.9
$ws.Range("M33").Value = -131.9

$ws.Range("H64").Value = 4496.5
$ws.Range("I64").Value = 4496.5
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 4496.5
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4248.5

$ws.Range("H67").Value = 4496.5
$ws.Range("I67").Value = 4496.5
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 4496.5
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -3638.5

$ws.Range("H70").Value = 4999.769
$ws.Range("J70").Value = 4999.769
$ws.Range("L70").Value = 14999.307
$ws.Range("N70").Value = -15539.307

$ws.Range("H73").Value = 4999.769
$ws.Range("J73").Value = 4999.769
$ws.Range("L73").Value = 14999.307
$ws.Range("N73").Value = -16871.307

$ws.Range("H132").Value = 21843.438
$ws.Range("I132").Value = 5791.3335
$ws.Range("K132").Value = 17374.0005
$ws.Range("M132").Value = -14844.0005

$ws.Range("H137").Value = 3470.6365
$ws.Range("I137").Value = 3703.4314
$ws.Range("J137").Value = 2679.1333
$ws.Range("K137").Value = 11110.2942
$ws.Range("L137").Value = 8037.3999
$ws.Range("M137").Value = -8560.2942
$ws.Range("N137").Value = -13137.3999

$ws.Range("H138").Value = 2881.8645
$ws.Range("J138").Value = 4829.6787
$ws.Range("L138").Value = 14489.0361
$ws.Range("N138").Value = -24769.0361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12992001
$ws.Range("I32").Value = 14290864
$ws.Range("K32").Value = 14290864
$ws.Range("M32").Value = -14290577

$ws.Range("H61").Value = 3088.1333
$ws.Range("I61").Value = 1666.2727
$ws.Range("J61").Value = 6998.25
$ws.Range("K61").Value = 1666.2727
$ws.Range("L61").Value = 6998.25
$ws.Range("M61").Value = -1454.2727
$ws.Range("N61").Value = -7422.25

$ws.Range("H74").Value = 3598
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 3598
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H132").Value = 1578.6666
$ws.Range("I132").Value = 1578.6666
$ws.Range("K132").Value = 4735.9998
$ws.Range("M132").Value = -2205.9998

$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws.Range("H136").Value = 3088.1333
$ws.Range("I136").Value = 1666.2727
$ws.Range("J136").Value = 6998.25
$ws.Range("K136").Value = 4998.8181
$ws.Range("L136").Value = 20994.75
$ws.Range("M136").Value = -2448.8181
$ws.Range("N136").Value = -26094.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7124.8
$ws.Range("I107").Value = 5781.125
$ws.Range("J107").Value = 12499.5
$ws.Range("K107").Value = 5781.125
$ws.Range("L107").Value = 12499.5
$ws.Range("M107").Value = -3861.125
$ws.Range("N107").Value = -16339.5

$ws.Range("H134").Value = 3198.7114
$ws.Range("I134").Value = 2880.3076
$ws.Range("J134").Value = 4153.923
$ws.Range("K134").Value = 8640.9228
$ws.Range("L134").Value = 12461.769
$ws.Range("M134").Value = -6105.9228
$ws.Range("N134").Value = -17531.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 784.3333
$ws.Range("I2").Value = 197
$ws.Range("J2").Value = 901.8
$ws.Range("K2").Value = 197
$ws.Range("L2").Value = 901.8
$ws.Range("M2").Value = -84
$ws.Range("N2").Value = -1127.8

$ws.Range("H31").Value = 2066.2727
$ws.Range("I31").Value = 1953.5555
$ws.Range("K31").Value = 1953.5555
$ws.Range("M31").Value = -1658.5555

$ws.Range("H34").Value = 2066.2727
$ws.Range("I34").Value = 1953.5555
$ws.Range("K34").Value = 1953.5555
$ws.Range("M34").Value = -1751.5555

$ws.Range("H62").Value = 12250
$ws.Range("I62").Value = 18000
$ws.Range("K62").Value = 18000
$ws.Range("M62").Value = -17376

$ws.Range("H65").Value = 12250
$ws.Range("I65").Value = 18000
$ws.Range("K65").Value = 90000
$ws.Range("M65").Value = -86880

$ws.Range("H107").Value = 12269.333
$ws.Range("I107").Value = 993.8
$ws.Range("J107").Value = 26363.75
$ws.Range("K107").Value = 993.8
$ws.Range("L107").Value = 26363.75
$ws.Range("M107").Value = 926.2
$ws.Range("N107").Value = -30203.75

$ws.Range("H132").Value = 2832.05
$ws.Range("I132").Value = 2861.9412
$ws.Range("K132").Value = 8585.8236
$ws.Range("M132").Value = -6055.8236

$ws.Range("H134").Value = 4157.067
$ws.Range("I134").Value = 3942.842
$ws.Range("J134").Value = 4527.091
$ws.Range("K134").Value = 11828.526
$ws.Range("L134").Value = 13581.273
$ws.Range("M134").Value = -9293.526
$ws.Range("N134").Value = -18651.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 167183.25
$ws.Range("J4").Value = 865.6667
$ws.Range("L4").Value = 2597.0001
$ws.Range("N4").Value = -2821.0001

$ws.Range("H113").Value = 1234
$ws.Range("J113").Value = 1128.5
$ws.Range("L113").Value = 3385.5
$ws.Range("N113").Value = -7725.5

$ws.Range("H139").Value = 166673840
$ws.Range("I139").Value = 500004000
$ws.Range("J139").Value = 8750
$ws.Range("K139").Value = 1500012000
$ws.Range("L139").Value = 26250
$ws.Range("M139").Value = -1500006860
$ws.Range("N139").Value = -36530

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 117.391304
$ws.Range("I2").Value = 78.933334
$ws.Range("J2").Value = 189.5
$ws.Range("K2").Value = 78.933334
$ws.Range("L2").Value = 189.5
$ws.Range("M2").Value = 34.066666
$ws.Range("N2").Value = -415.5

$ws.Range("H70").Value = 150394.88
$ws.Range("I70").Value = 229633.4
$ws.Range("K70").Value = 229633.4
$ws.Range("M70").Value = -229363.4

$ws.Range("H73").Value = 150394.88
$ws.Range("I73").Value = 229633.4
$ws.Range("K73").Value = 229633.4
$ws.Range("M73").Value = -228697.4

$ws.Range("H132").Value = 2806.7715
$ws.Range("I132").Value = 2784.5667
$ws.Range("K132").Value = 8353.7001
$ws.Range("M132").Value = -5823.7001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1420.409
$ws.Range("I55").Value = 966.6667
$ws.Range("K55").Value = 966.6667
$ws.Range("M55").Value = -793.6667

$ws.Range("H68").Value = 3972.923
$ws.Range("I68").Value = 3941.3333
$ws.Range("K68").Value = 3941.3333
$ws.Range("M68").Value = -3192.3333

$ws.Range("H71").Value = 3972.923
$ws.Range("I71").Value = 3941.3333
$ws.Range("K71").Value = 19706.6665
$ws.Range("M71").Value = -15962.6665

$ws.Range("H132").Value = 4574.4326
$ws.Range("I132").Value = 2692.3333
$ws.Range("K132").Value = 8076.999899999999
$ws.Range("M132").Value = -5546.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2528.4707
$ws.Range("I126").Value = 1073.375
$ws.Range("J126").Value = 3821.889
$ws.Range("K126").Value = 3220.125
$ws.Range("L126").Value = 11465.667
$ws.Range("M126").Value = -750.125
$ws.Range("N126").Value = -16405.667

$ws.Range("H132").Value = 3122.9546
$ws.Range("I132").Value = 2097.8
$ws.Range("K132").Value = 6293.400000000001
$ws.Range("M132").Value = -3763.400000000001

$ws.Range("H136").Value = 1621.1904
$ws.Range("I136").Value = 864.6667
$ws.Range("J136").Value = 3512.5
$ws.Range("K136").Value = 2594.0001
$ws.Range("L136").Value = 10537.5
$ws.Range("M136").Value = -44.0001000000002
$ws.Range("N136").Value = -15637.5
